$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 31, shifting existing rows 31-72 down to 32-73
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly data point
$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44671
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Magnum"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = 27000
$ws.Range("L31").Value = 28000
$ws.Range("M31").Value = 27500
$ws.Range("N31").Value = "$/saco 25 kilos"
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 1100
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
